$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-12 20:40:16"
$wsZhCn.Range("H3").Value = "2016-03-12 20:40:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-12 20:40:20"
$wsDeDe.Range("H3").Value = "2016-03-12 20:40:41"
